$d = $word.ActiveDocument

# --- Change 1: "In progress" -> "done" ------------------------------------
# The text lives in its own bold run ("Remove irrelevant attributes ... In
# progress"); a plain Find/Replace on the exact phrase only swaps the text
# inside that run and keeps the existing <w:b/><w:bCs/> formatting intact.
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("In progress", $false, $false, $false, $false, $false, $true, 1, $false, "done", 2)

# --- Change 2: append " done" (bold) after the "Review columns..." bullet -
# Grab the already-existing " done" formatted run (space + bold "done")
# that follows "Upload data" earlier in the list, so the new text reuses
# the exact same run split / bold formatting instead of merging into the
# preceding run.
$srcRange = $d.Content
$srcRange.Find.ClearFormatting()
$srcRange.Find.Execute("Upload data", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$srcStart = $srcRange.End
$srcSample = $d.Range($srcStart, $srcStart + 5)   # " done"

$formattedDone = $srcSample.FormattedText

$dstRange = $d.Content
$dstRange.Find.ClearFormatting()
$dstRange.Find.Execute("Review columns that have more than 50% n/a", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$dstEnd = $dstRange.End

$insertPoint = $d.Range($dstEnd, $dstEnd)
$insertPoint.FormattedText = $formattedDone
